$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place edits of shared-string runs) ---
# "Volume 31   Number  51" -> "...52"
$ws.Range("A8").Characters(21, 2).Text = "52"

# "Report Covering the Week  12/16/2024  Through  12/22/2024"
$ws.Range("C9").Characters(27, 10).Text = "12/23/2024"
$ws.Range("C9").Characters(48, 10).Text = "12/29/2024"

# --- Cells that flip between a numeric value and the "0" / "***.*" text placeholders ---
# These need both the cell style AND the underlying shared-string value changed, so we
# copy an already-styled donor cell (row 23, which is untouched by this edit) into place
# and then (for the two cells that become numeric) overwrite with the real number.

$ws.Range("C23").Copy($ws.Range("D15"))
$ws.Range("C23").Copy($ws.Range("D27"))
$ws.Range("C23").Copy($ws.Range("F27"))
$ws.Range("C23").Copy($ws.Range("D29"))
$ws.Range("C23").Copy($ws.Range("D30"))
$ws.Range("C23").Copy($ws.Range("D31"))
$ws.Range("E23").Copy($ws.Range("E15"))
$ws.Range("E23").Copy($ws.Range("E27"))
$ws.Range("E23").Copy($ws.Range("E29"))
$ws.Range("E23").Copy($ws.Range("E30"))
$ws.Range("E23").Copy($ws.Range("E31"))

$ws.Range("G23").Copy($ws.Range("D16"))
$ws.Range("D16").Value2 = 2
$ws.Range("H23").Copy($ws.Range("E16"))
$ws.Range("E16").Value2 = 50

# --- Remaining plain numeric value updates ---
$ws.Range("L15").Value2 = 23.076923076923
$ws.Range("C16").Value2 = 3
$ws.Range("F16").Value2 = 10
$ws.Range("H16").Value2 = 25
$ws.Range("I16").Value2 = 90
$ws.Range("J16").Value2 = 127
$ws.Range("K16").Value2 = -29.133858267716
$ws.Range("L16").Value2 = -42.675159235668
$ws.Range("M16").Value2 = -66.789667896679
$ws.Range("N16").Value2 = -94.642857142857
$ws.Range("C17").Value2 = 3
$ws.Range("D17").Value2 = 4
$ws.Range("E17").Value2 = -25
$ws.Range("G17").Value2 = 26
$ws.Range("H17").Value2 = -23.076923076923
$ws.Range("I17").Value2 = 270
$ws.Range("J17").Value2 = 360
$ws.Range("L17").Value2 = -21.511627906976
$ws.Range("M17").Value2 = -7.534246575342
$ws.Range("N17").Value2 = -70.779220779220
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 1
$ws.Range("E18").Value2 = 200
$ws.Range("F18").Value2 = 10
$ws.Range("G18").Value2 = 5
$ws.Range("H18").Value2 = 100
$ws.Range("I18").Value2 = 78
$ws.Range("J18").Value2 = 81
$ws.Range("K18").Value2 = -3.703703703703
$ws.Range("L18").Value2 = -42.222222222222
$ws.Range("M18").Value2 = -74.426229508196
$ws.Range("N18").Value2 = -95.822174611676
$ws.Range("C19").Value2 = 4
$ws.Range("D19").Value2 = 3
$ws.Range("E19").Value2 = 33.333333333333
$ws.Range("F19").Value2 = 8
$ws.Range("G19").Value2 = 23
$ws.Range("H19").Value2 = -65.217391304347
$ws.Range("I19").Value2 = 230
$ws.Range("J19").Value2 = 406
$ws.Range("K19").Value2 = -43.349753694581
$ws.Range("L19").Value2 = -39.790575916230
$ws.Range("M19").Value2 = -45.238095238095
$ws.Range("N19").Value2 = -66.521106259097
$ws.Range("D20").Value2 = 4
$ws.Range("F20").Value2 = 3
$ws.Range("G20").Value2 = 8
$ws.Range("H20").Value2 = -62.5
$ws.Range("J20").Value2 = 99
$ws.Range("K20").Value2 = -21.212121212121
$ws.Range("L20").Value2 = -46.575342465753
$ws.Range("N20").Value2 = -94.310722100656
$ws.Range("D21").Value2 = 14
$ws.Range("E21").Value2 = -7.142857142857
$ws.Range("F21").Value2 = 51
$ws.Range("G21").Value2 = 72
$ws.Range("H21").Value2 = -29.166666666666
$ws.Range("I21").Value2 = 767
$ws.Range("J21").Value2 = 1098
$ws.Range("K21").Value2 = -30.145719489981
$ws.Range("L21").Value2 = -35.328836424957
$ws.Range("M21").Value2 = -47.501711156742
$ws.Range("N21").Value2 = -88.447055279409
$ws.Range("L22").Value2 = -42.105263157894
$ws.Range("C24").Value2 = 26
$ws.Range("D24").Value2 = 10
$ws.Range("E24").Value2 = 160
$ws.Range("F24").Value2 = 80
$ws.Range("G24").Value2 = 80
$ws.Range("H24").Value2 = 0
$ws.Range("I24").Value2 = 912
$ws.Range("J24").Value2 = 1111
$ws.Range("K24").Value2 = -17.911791179117
$ws.Range("L24").Value2 = -17.540687160940
$ws.Range("M24").Value2 = -9.072781655034
$ws.Range("C25").Value2 = 6
$ws.Range("D25").Value2 = 3
$ws.Range("F25").Value2 = 13
$ws.Range("G25").Value2 = 16
$ws.Range("H25").Value2 = -18.75
$ws.Range("I25").Value2 = 156
$ws.Range("J25").Value2 = 239
$ws.Range("K25").Value2 = -34.728033472803
$ws.Range("L25").Value2 = -2.5
$ws.Range("C26").Value2 = 6
$ws.Range("D26").Value2 = 6
$ws.Range("E26").Value2 = 0
$ws.Range("F26").Value2 = 33
$ws.Range("G26").Value2 = 31
$ws.Range("H26").Value2 = 6.451612903225
$ws.Range("I26").Value2 = 521
$ws.Range("J26").Value2 = 503
$ws.Range("K26").Value2 = 3.578528827037
$ws.Range("L26").Value2 = -5.956678700361
$ws.Range("M26").Value2 = -37.153196622436
$ws.Range("H27").Value2 = -100
$ws.Range("L27").Value2 = -7.692307692307
$ws.Range("F28").Value2 = 1
$ws.Range("G28").Value2 = 5
$ws.Range("H28").Value2 = -80
$ws.Range("J28").Value2 = 58
$ws.Range("K28").Value2 = 10.344827586206
$ws.Range("M29").Value2 = -65.853658536585
$ws.Range("N29").Value2 = -86.274509803921
$ws.Range("M30").Value2 = -63.636363636363
$ws.Range("N30").Value2 = -87.878787878787
